$d = $word.ActiveDocument

# The underlying edit rotates the body text of several paragraphs while every
# paragraph keeps its own style/position/run-formatting exactly where it is:
#
#   "Objetivos" body paragraph              : A (Oferecer...)          -> C (Visão geral...)
#   "Docente(s) Responsável(eis)" bullet     : B (1033242 - Fábio...)   -> A (Oferecer...)
#   "Programa resumido" body paragraph       : C (Visão geral...)      -> D (long "- Introdução..." list)
#   "Programa" body paragraph                : D (long list)           -> E (Projetos pré-estruturados...)
#   "Avaliação" bullet - "Método:" run       : E (Projetos...)         -> F (As notas...)
#   "Avaliação" bullet - "Critério:" run     : F (As notas...)         -> G (Não há recuperação...)
#   "Avaliação" bullet - "Norma de..." run   : G (Não há recuperação...) -> H (bibliography)
#   "Bibliografia" body paragraph            : H (bibliography)         -> B (1033242 - Fábio...)
#
# The "D" block is a single run containing many <w:t>/<w:br/> pairs, so it is
# moved with a FormattedText range copy (preserving its internal structure)
# instead of a plain Find/Replace. The remaining (plain, single run) swaps are
# done through Find/Replace, routed through unique placeholder tokens first so
# that the rotated values never collide with one another while in flight.

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

$A = "Oferecer ao estudante uma visão ampla sobre aplicações contemporâneas e futuras de materiais poliméricos em áreas como a Engenharia de Materiais, Biomedicina, Nanotecnologia e outras."
$B = "1033242 - Fábio Herbst Florenzano"
$C = "Visão geral da Ciência de Polímeros; aplicações atuais e futuras de materiais poliméricos avançados; síntese de polímeros dirigida para aplicações específicas."
$E = "Projetos pré-estruturados de aplicações potenciais de polímeros realizados emgrupo com apresentação de texto escrito e seminário."
$F = "As notas (0 a 10) serão atribuídas conforme avaliação do projeto escrito e da apresentação do seminário. A nota mínima 5 é exigida para aprovação na disciplina."
$G = "Não há recuperação para esta disciplina."
$H = "1 - Akcelrud, L. Fundamentos da Ciência dos Polímeros. 1a. Edição. São Paulo:Editora Manole, 2006. ISBN: 978-8-85-2041561-0; 2 - AlMaadeed, M. A. A.; Ponnamma, D.; Carignano, M. A. Polymer Science and. 3 - Innovative Applications – Materials, Techniques and Future Developments, 1st edition. Amsterdã: Elsevier, 2020. ISBN: 978-0-12-816808-0. 4 - Canevarolo Jr., S. V.; Ciência dos Polímeros: Um Texto Básico para Tecnólogos e Engenheiros. 3a. Edição. São Paulo: ArtLiber, 2010. ISBN: 978-8-58-809810-7. 5 - Carraher, C. E. Introduction to Polymer Chemistry, 4th edition, Boca Raton: CRC Press – Taylor and Francis, 2017. ISBN: 978-1-4987-3761-6; 6 - Narain, R. Polymer Science and Nanotecnology: Fundamentals and Applications, 1 st edition. Amsterdã: Elsevier, 2020 ISBN: 978-0-12-816806-6; 7 - Lendlein, A. Shape-Memory Polymers, 1st Edition, Berlim: Springer-Velag, 2010. ISBN 978-3-642-12358-0; 8 - Rangelov, S.; Pispas, S. Polymer and Polymer-Hybrid Nanoparticles. 1st edition. 9 - Boca Raton: CRC Press, 2014. ISBN: 978-1-4398-6909-3"

# --- Step 1: tokenize every plain, single-run source text with a unique
#     placeholder BEFORE any literal text gets written anywhere, so later
#     writes can never be confused with text still to be located. ---
Replace-Text $A "@@SWAP_A@@"
Replace-Text $B "@@SWAP_B@@"
Replace-Text $E "@@SWAP_E@@"
Replace-Text $F "@@SWAP_F@@"
Replace-Text $G "@@SWAP_G@@"
Replace-Text $H "@@SWAP_H@@"

# --- Step 2: move the multi-run "D" (program list) block with a range copy,
#     from the "Programa" paragraph into the "Programa resumido" paragraph.
#     (Captured now; its own original text, "C", is simply overwritten so it
#     never needs to be located through Find.) ---
$pProgramaResumidoBody = $d.Paragraphs.Item(10)   # currently holds C
$pProgramaBody         = $d.Paragraphs.Item(12)   # currently holds D (the long list)

$dBlock = $pProgramaBody.Range.FormattedText

$rResumido = $d.Range($pProgramaResumidoBody.Range.Start, $pProgramaResumidoBody.Range.End - 1)
$rResumido.FormattedText = $dBlock

# Re-fetch the paragraph (earlier range edits can shift offsets) and overwrite
# it with E's literal text. This is safe now because the original "E" text
# (at the "Método:" run) was already swapped for "@@SWAP_E@@" in step 1.
$pProgramaBody2 = $d.Paragraphs.Item(12)
$rPrograma = $d.Range($pProgramaBody2.Range.Start, $pProgramaBody2.Range.End - 1)
$rPrograma.Text = $E

# --- Step 3: replace every placeholder with its final destination text. ---
Replace-Text "@@SWAP_A@@" $C   # "Objetivos" body: A -> C
Replace-Text "@@SWAP_B@@" $A   # "Docente(s)" bullet: B -> A
Replace-Text "@@SWAP_E@@" $F   # "Método:" run: E -> F
Replace-Text "@@SWAP_F@@" $G   # "Critério:" run: F -> G
Replace-Text "@@SWAP_G@@" $H   # "Norma de recuperação:" run: G -> H
Replace-Text "@@SWAP_H@@" $B   # "Bibliografia" body: H -> B
